$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.418769836425781
$ws.Range("B1").Value = 1.575741171836853
$ws.Range("C1").Value = 1.589781999588013
$ws.Range("D1").Value = 2.021628618240356
$ws.Range("E1").Value = 3.142244815826416
